$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Darshan"
$ws.Range("B3").Value = "Saman"
$ws.Range("B4").Value = "Shail"
$ws.Range("B5").Value = "Vaibhav"
$ws.Range("B6").Value = "John"

$ws.Range("C6").Select()
